# Weekly update: insert a new "Zanahoria" price record for Vega Monumental
# Concepción at row 135, pushing the existing rows 135-148 down to 136-149.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 135 (shifts old rows 135..148 -> 136..149,
# inheriting the D-column date style like Excel's native "Insert" command).
$ws.Rows(135).Insert()

# Populate the newly inserted row 135 with the new weekly record.
$ws.Range("A135").Value = 11
$ws.Range("B135").Value = "Vega Monumental Concepción"
$ws.Range("C135").Value = "Bíobío"
$ws.Range("D135").Value = 44516
$ws.Range("E135").Value = 8
$ws.Range("F135").Value = 100114013
$ws.Range("G135").Value = "Zanahoria"
$ws.Range("H135").Value = "Sin especificar"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 250
$ws.Range("K135").Value = 8000
$ws.Range("L135").Value = 8500
$ws.Range("M135").Value = 8300
$ws.Range("N135").Value = '$/saco 20 kilos'
$ws.Range("O135").Value = "Chillán"
$ws.Range("P135").Value = 415
$ws.Range("Q135").Value = 20
$ws.Range("R135").Value = "Hortaliza"
